$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the downloaded faostat bulk-download URLs, one per row, in the
# order they were fetched.
$ws.Range("A3").Value = "http://fenixservices.fao.org/faostat/static/bulkdownloads/Trade_Crops_Livestock_E_All_Data.zip"
$ws.Range("A4").Value = "http://fenixservices.fao.org/faostat/static/bulkdownloads/FoodSupply_Crops_E_All_Data.zip"
$ws.Range("A5").Value = "http://fenixservices.fao.org/faostat/static/bulkdownloads/Population_E_All_Data.zip"
$ws.Range("A6").Value = "http://fenixservices.fao.org/faostat/static/bulkdownloads/Value_of_Production_E_All_Data.zip"
$ws.Range("A7").Value = "http://fenixservices.fao.org/faostat/static/bulkdownloads/FoodBalanceSheets_E_All_Data.zip"

# Rows 5 and 6 were actually swapped when reviewing the list.
$ws.Range("A5").Value = "http://fenixservices.fao.org/faostat/static/bulkdownloads/Value_of_Production_E_All_Data.zip"
$ws.Range("A6").Value = "http://fenixservices.fao.org/faostat/static/bulkdownloads/Population_E_All_Data.zip"

$ws.Range("A7").Select()
